$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 50

# Columns A-D hold text values (date/time/weekday/week as strings).
# Force text format so Excel doesn't auto-convert them to a date serial
# number / numeric value, then clear the formatting again so the new
# cells end up with the default (unstyled) cell format, same as the
# other data rows.
$textRange = $ws.Range("A" + $row + ":D" + $row)
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-06-16"
$ws.Cells.Item($row, 2).Value = "17:50:00"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "24"

$textRange.ClearFormats()

# Columns E-T hold plain numeric values.
$ws.Cells.Item($row, 5).Value = 121806
$ws.Cells.Item($row, 6).Value = 132961
$ws.Cells.Item($row, 7).Value = 162001
$ws.Cells.Item($row, 8).Value = 133113
$ws.Cells.Item($row, 9).Value = 176942
$ws.Cells.Item($row, 10).Value = 114461
$ws.Cells.Item($row, 11).Value = 200634
$ws.Cells.Item($row, 12).Value = 224497
$ws.Cells.Item($row, 13).Value = 174697
$ws.Cells.Item($row, 14).Value = 103107
$ws.Cells.Item($row, 15).Value = 39125
$ws.Cells.Item($row, 16).Value = 33993
$ws.Cells.Item($row, 17).Value = 51726
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36748
$ws.Cells.Item($row, 20).Value = -1

$wb.Save()
